$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Snapshot the existing B (User Story) / C (API URL) columns for rows 6..33
#        BEFORE any writes, so the shift-down below doesn't clobber source data.
$oldB = @{}
$oldC = @{}
for ($r = 6; $r -le 33; $r++) {
    $oldB[$r] = $ws.Cells.Item($r, 2).Value2
    $oldC[$r] = $ws.Cells.Item($r, 3).Value2
}

# --- 2. Shift rows 6..33 down by one row in columns B/C only (column A - the
#        US ID - stays put per row, matching the target diff). Walk from the
#        bottom up so we never overwrite a value we still need to read.
for ($r = 33; $r -ge 6; $r--) {
    $srcRow = $r - 1
    $ws.Cells.Item($r, 2).Value2 = $oldB[$srcRow]
    if ($null -ne $oldC[$srcRow] -and $oldC[$srcRow] -ne "") {
        $ws.Cells.Item($r, 3).Value2 = $oldC[$srcRow]
    }
}

# --- 3. New row 34 gets the content that fell off the bottom (old row 33's
#        B value - "Employer pages ... ") plus a brand new US ID.
$ws.Cells.Item(34, 1).Value2 = "JB_US_33"
$ws.Cells.Item(34, 2).Value2 = $oldB[33]

# --- 4. Row 6 becomes the new "build resume" user story / API entry.
$ws.Cells.Item(6, 2).Value2 = "As a Candidate I am able to build the resume"
$ws.Cells.Item(6, 3).Value2 = "/api/resume/build/candidate"

# --- 5. Update the view: drop the scrolled-to topLeftCell and move the
#        active selection to C12 (matches the committed sheetView).
$ws.Application.Goto($ws.Range("A1"), $true)
$ws.Range("C12").Select()
